# Daily attendance processing - swap the first two comma-separated
# entries in the "Recorded By" column (G) for every data row that has
# two or more entries. Entries beyond the first two (if any) are left
# in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text.Split(",")

    if ($parts.Length -ge 2) {
        $p0 = $parts[0].Trim()
        $p1 = $parts[1].Trim()

        $rest = ""
        if ($parts.Length -gt 2) {
            for ($i = 2; $i -lt $parts.Length; $i++) {
                $rest = $rest + ", " + $parts[$i].Trim()
            }
        }

        $newVal = $p1 + ", " + $p0 + $rest
        $cell.Value = $newVal
    }
}
